$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($CellRef, $Val)
    $rng = $ws.Range($CellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $Val
    $rng.Style = "Normal"
}

$ws.Range("D2").Value = "62.037.88"
$ws.Range("E2").Value = "  -0.43%  "
$ws.Range("D3").Value = "3.007.79"
$ws.Range("E3").Value = "  +0.17%  "
Set-TextValue "D4" "1.00"
$ws.Range("E4").Value = "  +0.03%  "
Set-TextValue "D5" "597.15"
$ws.Range("E5").Value = "  +2.52%  "
Set-TextValue "D6" "146.98"
$ws.Range("E6").Value = "  +0.84%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "3.007.91"
$ws.Range("E8").Value = "  +0.25%  "
$ws.Range("E9").Value = "  -1.99%  "
$ws.Range("B10").Value = "Toncoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue "D10" "6.28"
$ws.Range("E10").Value = "  +8.04%  "
$ws.Range("B11").Value = "Dogecoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
Set-TextValue "D11" "0.149"
$ws.Range("E11").Value = "  +0.94%  "
$ws.Range("E12").Value = "  -0.50%  "
Set-TextValue "D13" "0.0000229"
$ws.Range("E13").Value = "  +0.46%  "
Set-TextValue "D14" "34.39"
$ws.Range("E14").Value = "  -0.17%  "
$ws.Range("E15").Value = "  +3.08%  "
$ws.Range("D16").Value = "3.503.94"
$ws.Range("E16").Value = "  +0.30%  "
$ws.Range("D17").Value = "61.970.00"
$ws.Range("E17").Value = "  -0.36%  "
$ws.Range("E18").Value = "  -2.25%  "
$ws.Range("D19").Value = "3.012.21"
$ws.Range("E19").Value = "  +0.31%  "
Set-TextValue "D20" "449.36"
$ws.Range("E20").Value = "  -2.42%  "
Set-TextValue "D21" "14.10"
$ws.Range("E21").Value = "  +0.94%  "
Set-TextValue "D22" "0.686"
$ws.Range("E22").Value = "  +0.08%  "
$ws.Range("E23").Value = "  -0.83%  "
Set-TextValue "D24" "81.80"
$ws.Range("E24").Value = "  +0.42%  "
Set-TextValue "D25" "11.10"
$ws.Range("E25").Value = "  +11.53%  "
Set-TextValue "D26" "2.25"
$ws.Range("E26").Value = "  +1.96%  "
Set-TextValue "D27" "12.18"
$ws.Range("E27").Value = "  -0.71%  "
$ws.Range("E28").Value = "  +0.18%  "
Set-TextValue "D29" "2.71"
$ws.Range("E29").Value = "  +3.69%  "
Set-TextValue "D30" "1.00"
$ws.Range("E30").Value = "  +0.25%  "
Set-TextValue "D31" "7.21"
$ws.Range("E31").Value = "  +2.67%  "
Set-TextValue "D32" "2.09"
$ws.Range("E32").Value = "  +0.36%  "
Set-TextValue "D33" "27.33"
$ws.Range("E33").Value = "  -3.08%  "
Set-TextValue "D34" "0.110"
$ws.Range("E34").Value = "  +2.07%  "
$ws.Range("D35").Value = "0.0₃0836"
$ws.Range("E35").Value = "  +5.20%  "
Set-TextValue "D36" "1.02"
$ws.Range("E36").Value = "  -0.20%  "
Set-TextValue "D37" "5.80"
$ws.Range("E37").Value = "  +0.81%  "
Set-TextValue "D38" "50.52"
$ws.Range("E38").Value = "  +0.54%  "
Set-TextValue "D39" "2.05"
$ws.Range("E39").Value = "  -2.15%  "
Set-TextValue "D40" "9.00"
$ws.Range("E40").Value = "  -2.20%  "
$ws.Range("E41").Value = "  +8.36%  "
Set-TextValue "D42" "2.91"
$ws.Range("E42").Value = "  +1.16%  "
Set-TextValue "D43" "398.94"
$ws.Range("E43").Value = "  +1.38%  "
Set-TextValue "D44" "40.76"
$ws.Range("E44").Value = "  +11.30%  "
$ws.Range("E45").Value = "  +0.90%  "
Set-TextValue "D46" "0.0352"
$ws.Range("E46").Value = "  -1.53%  "
$ws.Range("D47").Value = "2.714.65"
$ws.Range("E47").Value = "  -0.26%  "
Set-TextValue "D48" "132.39"
$ws.Range("E48").Value = "  +3.09%  "
$ws.Range("E50").Value = "  -0.28%  "
$ws.Range("E51").Value = "  -1.51%  "
